# Update the quiz word list: scores bumped up and rows re-sorted by the
# new score (ties keep the app's ordering), giving the sheet a more
# symmetric look as the rows grow toward the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("miettiä",         "to think (over), to consider", 1),
    @("jossa",            "in which",                      1),
    @("jossain",          "somewhere",                     2),
    @("jotka",            "who/that/which (relative)",     2),
    @("olla ostoksilla",  "to be shopping",                3),
    @("joku",             "someone",                       3),
    @("tyhjä",            "empty",                         3),
    @("yhtään",           "any (at all)",                  3),
    @("virhe",            "mistake, error",                3),
    @("päivällinen",      "dinner",                        3),
    @("timantti",         "diamond",                       3),
    @("housut",           "pants",                         3),
    @("mekko",            "dress",                         3),
    @("kenkä",            "shoe",                          3),
    @("näytelmä",         "play (theatre)",                3),
    @("olut",             "beer",                          3),
    @("keskusta",         "city center",                   4),
    @("asukas",           "inhabitant",                    4),
    @("keittiö",          "kitchen",                       4),
    @("ruuhka",           "traffic jam",                   4),
    @("annos",            "portion, dish",                 4),
    @("auttaa",           "to help",                       4),
    @("yliopisto",        "university",                    4)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $word = $rows[$i][0]
    $meaning = $rows[$i][1]
    $score = $rows[$i][2]

    $ws.Cells.Item($r, 1).Value = $word
    $ws.Cells.Item($r, 2).Value = $meaning
    $ws.Cells.Item($r, 3).Value = $score
}
